$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.127881588408715
$ws.Range("C2").Value = 10.29869402782916
$ws.Range("D2").Value = 337.1190423067083
$ws.Range("E2").Value = 616238.5361209477
$ws.Range("G2").Value = 616586.0817388706
